# Slide 5 ("本日説明する内容"), content placeholder shape (shape index 2).
# Item 1: "... (田浦; 20分)"  -> "... (田浦; 25分)"   (time allocation 20 -> 25)
# Item 2: "... (田浦; 20分)"  -> "... (田浦; 15分)"   (time allocation 20 -> 15)
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Item 1 ("1. 授業に必要なICTシステムの概要 (田浦; 20分)"): change the minutes "20" -> "25".
$tr.Characters(26, 2).Text = "25"

# Item 2 ("2. 2020年度振り返り (田浦; 20分)"): change the minutes "20" -> "15"
# (the "; " and "20" runs collapse into a single "; 15" run).
$tr.Characters(49, 4).Text = "; 15"
